$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Validator" sheet right after "Orchestrator" ---
# Clone the "Orchestrator" sheet so the new sheet inherits the same sheet
# formatting (outline props, page setup, margins, etc.) as its siblings,
# then rename it and replace its contents with the Validator header row.
$orchestrator = $wb.Worksheets.Item("Orchestrator")
$orchestrator.Copy($null, $orchestrator)
$validator = $wb.Worksheets.Item("Orchestrator (2)")
$validator.Name = "Validator"

$validator.Cells.Clear()
$validator.Range("A1").Value = "name"
$validator.Range("B1").Value = "url"
$validator.Range("C1").Value = "schema_url"

# --- 2. Update the "Benchmark" sheet header row with the new fields ---
# (version, storage and validator are new columns; the rest shift right)
$benchmark = $wb.Worksheets.Item("Benchmark")

$benchmark.Range("A1").Value = "version"
$benchmark.Range("B1").Value = "platform"
$benchmark.Range("C1").Value = "storage"
$benchmark.Range("D1").Value = "orchestrator"
$benchmark.Range("E1").Value = "validator"
$benchmark.Range("F1").Value = "steps"
$benchmark.Range("G1").Value = "id"
$benchmark.Range("H1").Value = "name"
$benchmark.Range("I1").Value = "description"
